$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D. This shifts the existing D:K data
# (and its formatting) one column to the right, into E:L.
$ws.Columns("D:D").Insert()

# Copy the number formatting (date format in row 7/38/80, thousands
# format elsewhere) from the old "D" data (now shifted to E) into the
# freshly inserted, still-blank column D.
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Rows 36/37 and 78/79 have no data in D:K (37/79 are section headers,
# 36/78 are fully blank spacer rows that don't even exist in the
# sheetData) - the paste above spuriously stamped a blank formatted
# cell there, remove it so those rows stay exactly as they were.
$ws.Range("D36").Clear()
$ws.Range("D37").Clear()
$ws.Range("D78").Clear()
$ws.Range("D79").Clear()

# Populate the new column D with the newest reporting period's figures.
$ws.Cells.Item(7, 4).Value2 = 43465
$ws.Cells.Item(8, 4).Value2 = 1400000
$ws.Cells.Item(9, 4).Value2 = 830000
$ws.Cells.Item(10, 4).Value2 = 570000
$ws.Cells.Item(12, 4).Value2 = 52900
$ws.Cells.Item(13, 4).Value2 = 0
$ws.Cells.Item(14, 4).Value2 = 147000
$ws.Cells.Item(15, 4).Value2 = 0
$ws.Cells.Item(17, 4).Value2 = 1378400
$ws.Cells.Item(18, 4).Value2 = 21600
$ws.Cells.Item(20, 4).Value2 = 1500
$ws.Cells.Item(21, 4).Value2 = 97800
$ws.Cells.Item(22, 4).Value2 = 700
$ws.Cells.Item(23, 4).Value2 = 22400
$ws.Cells.Item(24, 4).Value2 = -7200
$ws.Cells.Item(25, 4).Value2 = 0
$ws.Cells.Item(26, 4).Value2 = 29600
$ws.Cells.Item(27, 4).Value2 = 29600
$ws.Cells.Item(28, 4).Value2 = 0
$ws.Cells.Item(29, 4).Value2 = -800
$ws.Cells.Item(30, 4).Value2 = 0
$ws.Cells.Item(31, 4).Value2 = 0
$ws.Cells.Item(32, 4).Value2 = -1500
$ws.Cells.Item(33, 4).Value2 = 28800
$ws.Cells.Item(34, 4).Value2 = 0
$ws.Cells.Item(35, 4).Value2 = 28800
$ws.Cells.Item(38, 4).Value2 = 43465
$ws.Cells.Item(41, 4).Value2 = 344800
$ws.Cells.Item(42, 4).Value2 = 37300
$ws.Cells.Item(43, 4).Value2 = 212900
$ws.Cells.Item(44, 4).Value2 = 311200
$ws.Cells.Item(45, 4).Value2 = 20900
$ws.Cells.Item(46, 4).Value2 = 927000
$ws.Cells.Item(47, 4).Value2 = 2000
$ws.Cells.Item(48, 4).Value2 = 432600
$ws.Cells.Item(49, 4).Value2 = 144600
$ws.Cells.Item(50, 4).Value2 = 0
$ws.Cells.Item(51, 4).Value2 = 0
$ws.Cells.Item(52, 4).Value2 = 79100
$ws.Cells.Item(53, 4).Value2 = 0
$ws.Cells.Item(54, 4).Value2 = 1585400
$ws.Cells.Item(57, 4).Value2 = 120500
$ws.Cells.Item(58, 4).Value2 = 0
$ws.Cells.Item(59, 4).Value2 = 128800
$ws.Cells.Item(60, 4).Value2 = 249300
$ws.Cells.Item(61, 4).Value2 = 0
$ws.Cells.Item(62, 4).Value2 = 72400
$ws.Cells.Item(63, 4).Value2 = 0
$ws.Cells.Item(64, 4).Value2 = 0
$ws.Cells.Item(65, 4).Value2 = 0
$ws.Cells.Item(66, 4).Value2 = 321700
$ws.Cells.Item(68, 4).Value2 = 0
$ws.Cells.Item(69, 4).Value2 = 0
$ws.Cells.Item(70, 4).Value2 = 0
$ws.Cells.Item(71, 4).Value2 = 0
$ws.Cells.Item(72, 4).Value2 = 620700
$ws.Cells.Item(73, 4).Value2 = 0
$ws.Cells.Item(74, 4).Value2 = 0
$ws.Cells.Item(75, 4).Value2 = 0
$ws.Cells.Item(76, 4).Value2 = 1263700
$ws.Cells.Item(77, 4).Value2 = 0
$ws.Cells.Item(80, 4).Value2 = 43465
$ws.Cells.Item(81, 4).Value2 = 28800
$ws.Cells.Item(83, 4).Value2 = 74700
$ws.Cells.Item(84, 4).Value2 = 0
$ws.Cells.Item(85, 4).Value2 = 0
$ws.Cells.Item(86, 4).Value2 = 0
$ws.Cells.Item(87, 4).Value2 = 0
$ws.Cells.Item(88, 4).Value2 = 0
$ws.Cells.Item(89, 4).Value2 = 160200
$ws.Cells.Item(92, 4).Value2 = 0
$ws.Cells.Item(93, 4).Value2 = 0
$ws.Cells.Item(94, 4).Value2 = -103400
$ws.Cells.Item(96, 4).Value2 = 0
$ws.Cells.Item(97, 4).Value2 = 0
$ws.Cells.Item(98, 4).Value2 = 0
$ws.Cells.Item(99, 4).Value2 = 0
$ws.Cells.Item(100, 4).Value2 = 8000
$ws.Cells.Item(101, 4).Value2 = -10200
$ws.Cells.Item(102, 4).Value2 = 54700

# Row 91 ("Changes In Accounts Receivables") also received corrected
# historical figures alongside the new column, not just a pure shift.
$ws.Cells.Item(91, 4).Value2 = -92700
$ws.Cells.Item(91, 5).Value2 = -74500
$ws.Cells.Item(91, 6).Value2 = -23400
$ws.Cells.Item(91, 7).Value2 = -13000
$ws.Cells.Item(91, 8).Value2 = -16600
$ws.Cells.Item(91, 9).Value2 = -18400
$ws.Cells.Item(91, 10).Value2 = -19200
$ws.Cells.Item(91, 11).Value2 = -15800

Write-Output "done"
